# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2410
#   *_new  -> *_FV2504
# Also freeze the header row and turn the data range into a proper Excel
# Table (ListObject) with an AutoFilter, matching the reference workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row ------------------------------------------------
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = ($oldHeaders[$i] -replace "_old$", "_FV2410")
}

# column 11 is "diff" and is left untouched

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = ($newHeaders[$i] -replace "_new$", "_FV2504")
}

# --- 2. Freeze the header row -------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table ------------------------------
$dataRange = $ws.Range("A1:U60")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""
